$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New CPF (C), RG (D), PIS/NIS (G) values per row (rows 2-36)
$values = @{
    2  = @{ C = "99123456789"; D = "7392846"; G = "83920471562" }
    3  = @{ C = "98765432109"; D = "5918374"; G = "74829103647" }
    4  = @{ C = "93456789012"; D = "8203947"; G = "98374612058" }
    5  = @{ C = "94567890123"; D = "6739201"; G = "75938462019" }
    6  = @{ C = "94321098765"; D = "7029381"; G = "61092837456" }
    7  = @{ C = "98765432098"; D = "6482917"; G = "79482013647" }
    8  = @{ C = "94567892345"; D = "7193846"; G = "70293847156" }
    9  = @{ C = "92345679876"; D = "8039172"; G = "64829103647" }
    10 = @{ C = "96543219876"; D = "7659203"; G = "71938462019" }
    11 = @{ C = "95678431234"; D = "8491736"; G = "80391745612" }
    12 = @{ C = "94567893456"; D = "7291836"; G = "71092837456" }
    13 = @{ C = "92345678098"; D = "7902847"; G = "68492013647" }
    14 = @{ C = "96543218765"; D = "7819204"; G = "72918347156" }
    15 = @{ C = "95678430123"; D = "6038471"; G = "79028413647" }
    16 = @{ C = "98765430986"; D = "6891735"; G = "78192047156" }
    17 = @{ C = "92345677098"; D = "6203948"; G = "68917362019" }
    18 = @{ C = "96543217654"; D = "6309281"; G = "75829103647" }
    19 = @{ C = "95678429012"; D = "7593841"; G = "62039447156" }
    20 = @{ C = "94321095432"; D = "8129374"; G = "79928462019" }
    21 = @{ C = "96543216543"; D = "6891734"; G = "62938462019" }
    22 = @{ C = "95678428098"; D = "7938472"; G = "69582047156" }
    23 = @{ C = "93456785432"; D = "6829174"; G = "79384762019" }
    24 = @{ C = "97890128901"; D = "6739202"; G = "68291713647" }
    25 = @{ C = "96543215432"; D = "7593844"; G = "74829162019" }
    26 = @{ C = "94567897890"; D = "6482916"; G = "79482047155" }
    27 = @{ C = "92345674098"; D = "5839203"; G = "83917413647" }
    28 = @{ C = "96543214321"; D = "8039171"; G = "64829147156" }
    29 = @{ C = "94321089876"; D = "6891733"; G = "62938413647" }
    30 = @{ C = "98765423219"; D = "7938471"; G = "69582062019" }
    31 = @{ C = "94567891234"; D = "7593846"; G = "68917313647" }
    32 = @{ C = "92345670098"; D = "6829173"; G = "79384747156" }
    33 = @{ C = "97890133456"; D = "7482914"; G = "61092813647" }
    34 = @{ C = "94321088765"; D = "7593847"; G = "74829113647" }
    35 = @{ C = "94567892345"; D = "8391744"; G = "75938462016" }
    36 = @{ C = "92345669098"; D = "7029383"; G = "61092847155" }
}

foreach ($row in $values.Keys) {
    $rowData = $values[$row]
    $ws.Cells.Item($row, 3).Value = [double]$rowData.C
    $ws.Cells.Item($row, 4).Value = [double]$rowData.D
    $ws.Cells.Item($row, 7).Value = [double]$rowData.G
}

# Update the "Sigla" labels used in column J (shared strings reused across rows)
$excel.Cells.Replace("SISFIN - SP", "SIGLA 5", 1) | Out-Null
$excel.Cells.Replace("SIAPI", "SIGLA 1", 1) | Out-Null
$excel.Cells.Replace("SIPEN", "SIGLA 2", 1) | Out-Null
$excel.Cells.Replace("SISFIN - DF", "SIGLA 4", 1) | Out-Null
